# Apply edit: add "happy number" row (with hash) to the 哈希 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("哈希")
$ws.Activate()

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 202
$ws.Cells.Item(4, 3).Value = "编写一个算法来判断一个数 n 是不是快乐数。 `n[快乐数」定义为：对于一个正整数，每一次将该数替换为它每个位置上的数字的平方和，然后重复这个过程直到这个数变为 1，也可能是 无限循环 但始终变不到 1。`n如果 可以变为 1，那么这个数就是快乐数。 "
$ws.Cells.Item(4, 4).Value = "0 这个数字是否在集合中出现过，出现过就false，没出现过就加入集合，继续1`n1 拆解一个整数，分别将各位相加`n2 判断这个数的各位平方和是否是1`n3 是1，就返回true`n4 不是1，将平方之和作为新的数字`n5 继续循环"
$ws.Cells.Item(4, 5).Value = "集合`n整数拆解"
$ws.Cells.Item(4, 6).Value = "未知"
$ws.Cells.Item(4, 7).Value = "未知"

$ws.Rows.Item(4).RowHeight = 220

$ws.Range("D4").Select()
